$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Row 2 - ANATOMY session 1: reorder "Recorded By" list
$ws.Range("G2").Value = "System, Amira.Sobhy@med.asu.edu.eg, servinaz@med.asu.edu.eg, gehanadel@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"

# Row 9 - HISTOLOGY session 1: reorder "Recorded By" list
$ws.Range("G9").Value = "Shimaa.ashraf@med.asu.edu.eg, Safa.hany@med.asu.edu.eg"

# Row 10 - HISTOLOGY Average Attendance % statistic update
# (leading apostrophe keeps this a literal text value "21.6%" instead of
#  Excel auto-converting the percent-looking string into a numeric 0.216)
$ws.Range("L10").Value = "'21.6%"

# Row 14 - PARASITOLOGY session 1: Students count update
$ws.Range("H14").Value = "32/251"

# Row 15 - PARASITOLOGY Avg Attendance % statistic update
$ws.Range("S15").Value = "'21.6%"

# Row 28 - PHYSIOLOGY session 1: reorder "Recorded By" list
$ws.Range("G28").Value = "maryam.ashraf@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg"
